$wb = $excel.ActiveWorkbook

# Row 9 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 248.8421
$ws.Cells.Item(9, 9).Value = 228.86667
$ws.Cells.Item(9, 11).Value = 228.86667
$ws.Cells.Item(9, 13).Value = -59.86667

# Row 29 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 4148.8335
$ws.Cells.Item(29, 10).Value = 4663.5835
$ws.Cells.Item(29, 12).Value = 13990.7505
$ws.Cells.Item(29, 14).Value = -14552.7505

# Row 33 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 300.5
$ws.Cells.Item(33, 9).Value = 201
$ws.Cells.Item(33, 11).Value = 201
$ws.Cells.Item(33, 13).Value = 28

# Row 40 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 6302.241
$ws.Cells.Item(40, 10).Value = 7951
$ws.Cells.Item(40, 12).Value = 7951
$ws.Cells.Item(40, 14).Value = -8301

# Row 106 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 425
$ws.Cells.Item(106, 9).Value = 425
$ws.Cells.Item(106, 11).Value = 425
$ws.Cells.Item(106, 13).Value = 206

# Row 132 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 31829.715
$ws.Cells.Item(132, 9).Value = 53760.75
$ws.Cells.Item(132, 11).Value = 161282.25
$ws.Cells.Item(132, 13).Value = -158752.25

# Row 138 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4690.4443
$ws.Cells.Item(138, 9).Value = 4598.5
$ws.Cells.Item(138, 10).Value = 4716.7144
$ws.Cells.Item(138, 11).Value = 13795.5
$ws.Cells.Item(138, 12).Value = 14150.1432
$ws.Cells.Item(138, 13).Value = -8655.5
$ws.Cells.Item(138, 14).Value = -24430.1432

# Row 32 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2176329.2
$ws.Cells.Item(32, 9).Value = 2250.0513
$ws.Cells.Item(32, 10).Value = 14289056
$ws.Cells.Item(32, 11).Value = 2250.0513
$ws.Cells.Item(32, 12).Value = 14289056
$ws.Cells.Item(32, 13).Value = -1963.0513
$ws.Cells.Item(32, 14).Value = -14289630

# Row 44 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 11772
$ws.Cells.Item(44, 10).Value = 11772
$ws.Cells.Item(44, 12).Value = 11772
$ws.Cells.Item(44, 14).Value = -12748

# Row 86 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).ClearContents()

# Row 88 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 899.5
$ws.Cells.Item(88, 9).Value = 800
$ws.Cells.Item(88, 11).Value = 800
$ws.Cells.Item(88, 13).Value = -394

# Row 89 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).ClearContents()

# Row 91 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 899.5
$ws.Cells.Item(91, 9).Value = 800
$ws.Cells.Item(91, 11).Value = 800
$ws.Cells.Item(91, 13).Value = 604

# Row 92 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 137499.75
$ws.Cells.Item(92, 10).Value = 137499.75
$ws.Cells.Item(92, 12).Value = 137499.75
$ws.Cells.Item(92, 14).Value = -142491.75

# Row 16 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 703
$ws.Cells.Item(16, 9).Value = 703
$ws.Cells.Item(16, 11).Value = 703
$ws.Cells.Item(16, 13).Value = -416

# Row 23 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 1999
$ws.Cells.Item(23, 9).Value = 1999
$ws.Cells.Item(23, 11).Value = 1999
$ws.Cells.Item(23, 13).Value = -1759

# Row 27 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(27, 8).Value = 1999
$ws.Cells.Item(27, 9).Value = 1999
$ws.Cells.Item(27, 11).Value = 1999
$ws.Cells.Item(27, 13).Value = -1807

# Row 35 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 1924.4445
$ws.Cells.Item(35, 9).Value = 936.6667
$ws.Cells.Item(35, 11).Value = 936.6667
$ws.Cells.Item(35, 13).Value = -642.6667

# Row 69 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(69, 8).Value = 39333
$ws.Cells.Item(69, 9).Value = 38999.5
$ws.Cells.Item(69, 11).Value = 38999.5
$ws.Cells.Item(69, 13).Value = -38250.5

# Row 72 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(72, 8).Value = 39333
$ws.Cells.Item(72, 9).Value = 38999.5
$ws.Cells.Item(72, 11).Value = 116998.5
$ws.Cells.Item(72, 13).Value = -113254.5

# Row 86 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 8649.75
$ws.Cells.Item(86, 9).Value = 8600
$ws.Cells.Item(86, 10).Value = 8666.333000000001
$ws.Cells.Item(86, 11).Value = 8600
$ws.Cells.Item(86, 12).Value = 8666.333000000001
$ws.Cells.Item(86, 13).Value = -7477
$ws.Cells.Item(86, 14).Value = -10912.333

# Row 89 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 8649.75
$ws.Cells.Item(89, 9).Value = 8600
$ws.Cells.Item(89, 10).Value = 8666.333000000001
$ws.Cells.Item(89, 11).Value = 43000
$ws.Cells.Item(89, 12).Value = 43331.665
$ws.Cells.Item(89, 13).Value = -37384
$ws.Cells.Item(89, 14).Value = -54563.665

# Row 113 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 703
$ws.Cells.Item(113, 9).Value = 703
$ws.Cells.Item(113, 11).Value = 703
$ws.Cells.Item(113, 13).Value = 1467

# Row 132 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1845.6
$ws.Cells.Item(132, 9).Value = 1247
$ws.Cells.Item(132, 11).Value = 3741
$ws.Cells.Item(132, 13).Value = -1211

# Row 134 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1106.9584
$ws.Cells.Item(134, 9).Value = 1106.9584
$ws.Cells.Item(134, 11).Value = 3320.8752
$ws.Cells.Item(134, 13).Value = -785.8751999999999

# Row 94 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 7427
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 13).ClearContents()

# Row 63 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 41332
$ws.Cells.Item(63, 10).Value = 41332
$ws.Cells.Item(63, 12).Value = 41332
$ws.Cells.Item(63, 14).Value = -42704

# Row 66 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(66, 8).Value = 41332
$ws.Cells.Item(66, 10).Value = 41332
$ws.Cells.Item(66, 12).Value = 123996
$ws.Cells.Item(66, 14).Value = -130860

# Row 92 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 5750.1113
$ws.Cells.Item(92, 10).Value = 5750.1113
$ws.Cells.Item(92, 12).Value = 5750.1113
$ws.Cells.Item(92, 14).Value = -9494.1113

# Row 113 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3805.7856
$ws.Cells.Item(113, 9).Value = 2116.4546
$ws.Cells.Item(113, 11).Value = 2116.4546
$ws.Cells.Item(113, 13).Value = 53.54539999999997

# Row 122 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 557776.3
$ws.Cells.Item(122, 9).Value = 1001100
$ws.Cells.Item(122, 10).Value = 3621.75
$ws.Cells.Item(122, 11).Value = 3003300
$ws.Cells.Item(122, 12).Value = 10865.25
$ws.Cells.Item(122, 13).Value = -3000850
$ws.Cells.Item(122, 14).Value = -15765.25

# Row 22 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3366.6667
$ws.Cells.Item(22, 10).Value = 4000
$ws.Cells.Item(22, 12).Value = 4000
$ws.Cells.Item(22, 14).Value = -4590

# Row 27 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3366.6667
$ws.Cells.Item(27, 10).Value = 4000
$ws.Cells.Item(27, 12).Value = 4000
$ws.Cells.Item(27, 14).Value = -4214

# Row 40 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4964
$ws.Cells.Item(40, 9).Value = 4946
$ws.Cells.Item(40, 11).Value = 4946
$ws.Cells.Item(40, 13).Value = -4810

# Row 46 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 5581.5454
$ws.Cells.Item(46, 9).Value = 2441.1667
$ws.Cells.Item(46, 10).Value = 9350
$ws.Cells.Item(46, 11).Value = 2441.1667
$ws.Cells.Item(46, 12).Value = 9350
$ws.Cells.Item(46, 13).Value = -2253.1667
$ws.Cells.Item(46, 14).Value = -9726

# Row 55 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1163
$ws.Cells.Item(55, 9).Value = 1123.125
$ws.Cells.Item(55, 11).Value = 1123.125
$ws.Cells.Item(55, 13).Value = -950.125

# Row 56 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(56, 8).Value = 16854.834
$ws.Cells.Item(56, 9).Value = 7971.6665
$ws.Cells.Item(56, 10).Value = 25738
$ws.Cells.Item(56, 11).Value = 7971.6665
$ws.Cells.Item(56, 12).Value = 25738
$ws.Cells.Item(56, 13).Value = -7280.6665
$ws.Cells.Item(56, 14).Value = -27120

# Row 74 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(74, 8).Value = 20000
$ws.Cells.Item(74, 9).Value = 20000
$ws.Cells.Item(74, 11).Value = 20000
$ws.Cells.Item(74, 13).Value = -19002

# Row 77 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(77, 8).Value = 20000
$ws.Cells.Item(77, 9).Value = 20000
$ws.Cells.Item(77, 11).Value = 60000
$ws.Cells.Item(77, 13).Value = -55008

# Row 93 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 8498.777
$ws.Cells.Item(93, 9).Value = 12331.667
$ws.Cells.Item(93, 10).Value = 6582.3335
$ws.Cells.Item(93, 11).Value = 12331.667
$ws.Cells.Item(93, 12).Value = 6582.3335
$ws.Cells.Item(93, 13).Value = -11083.667
$ws.Cells.Item(93, 14).Value = -9078.333500000001

# Row 100 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 9263
$ws.Cells.Item(100, 10).Value = 9999.833000000001
$ws.Cells.Item(100, 12).Value = 9999.833000000001
$ws.Cells.Item(100, 14).Value = -11081.833

# Row 51 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 10000
$ws.Cells.Item(51, 9).Value = 10000
$ws.Cells.Item(51, 11).Value = 10000
$ws.Cells.Item(51, 13).Value = -9490

# Row 52 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 39999
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()

# Row 107 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2268.8572
$ws.Cells.Item(107, 9).Value = 2658.6
$ws.Cells.Item(107, 11).Value = 7975.799999999999
$ws.Cells.Item(107, 13).Value = -6055.799999999999

# Row 112 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(112, 8).Value = 38000
$ws.Cells.Item(112, 10).Value = 38000
$ws.Cells.Item(112, 12).Value = 38000
$ws.Cells.Item(112, 14).Value = -40954
